$d = $word.ActiveDocument

$ops = @(
    ,@("Microbiology", "Microbiology: from theory to practice")
    ,@("Créditos-aula: 3`v", "Créditos-aula: 6`v")
    ,@("Carga horária: 45 h`v", "Carga horária: 90 h`v")
    ,@("Ativação: 01/01/2025`v", "Ativação: 15/07/2025`v")
    ,@("Fornecer aos alunos do curso de Engenharia Bioquímica, conhecimentos basicos de microbiologia e sua importância na agricultura, meio ambiente, medicina e indústria.", "Aspectos teóricos (3 créditos): Conceitos gerais de microbiologia; Filogenia Microbiana; Caracterização dos principais grupos; Cultivo e crescimento microbiano; Genética de Microrganismos; Ecologia e Interações Microbianas; Vírus e fagos; `v`vPráticas em laboratório (3 créditos): Técnicas de microscopia para visualização de células microbianas. Isolamento, cultivo e controle do crescimento de microrganismos. Quantificação de células. Identificação de microrganismos. Análise microbiológica de diferentes tipos de amostras.")
    ,@("Provide students in the Biochemical Engineering course with basic knowledge of microbiology and its importance in agriculture, the environment, medicine, and industry.", "Theoretical Aspects (3 credits): General microbiology concepts; Microbial Phylogeny; Characterization of main groups; Microbial cultivation and growth; Microorganism Genetics; Microbial Ecology and Interactions; Viruses and Phages.`v`vLaboratory Practices (3 credits): Microscopy techniques for visualizing microbial cells. Isolation, cultivation, and growth control of microorganisms. Cell quantification. Microorganism identification. Microbiological analysis of different sample types")
    ,@("4873328 - Fernando Segato`v", "Proporcionar aos alunos conhecimentos fundamentais sobre os principais grupos de microrganismos, abordando sua fisiologia, genética e funções em processos biotecnológicos, bem como desenvolver habilidades práticas em técnicas microbiológicas para isolamento, cultivo e caracterização desses organismos. Adicionalmente, incentivar a interpretação crítica de resultados experimentais.`v")
    ,@("8853480 - Tatiane da Franca Silva", "Aulas teóricas: `v`v1-Histórico da microbiologia, importância industrial dos microrganismos, taxonomia molecular e filogenia, fisiologia e caracterização dos diferentes grupos: bactérias, fungos, vírus e arqueas;`v2- Nutrição microbiana; meios de cultura; fatores ambientais; reprodução e crescimento microbiano; medidas de controle de crescimento; métodos de quantificação de microrganismos;`v3-Organização do genoma microbiano. Transferência genética: Conjugação, transformação e transdução. Regulação da expressão gênica: sistema Operon;`v4-Papel ecológico dos microrganismos em biociclos. Bio indicadores; Relações simbióticas;`v5-Genoma viral e sua diversidade;`v `vAulas práticas: `v`v1.Microscopia ótica e técnicas de coloração;`v2.Preparo e esterilização de meios de cultura - exigências nutricionais;`v3.Técnicas de semeadura de microrganismos e isolamento de linhagem em cultura pura;`v4.Ação de agentes físicos e químicos sobre o crescimento microbiano;`v5.Técnicas de quantificação de microrganismos: Unidade formadora de colônias (UFC), Turbidimetria e Câmara de contagem celular em microscópio`v6.Testes Bioquímicos na identificação e caracterização de microorganismos`v7.Análise bacteriana de amostras de água e leite")
    ,@("Histórico da microbiologia, microbiologia industrial, filogênia microbiana, caracterização dos microrganismos, nutrição e cultivo de microrganismos, virus, fungos filamentosos, leveduras, micro-algas, bactérias.", "A avaliação será composta por provas, exercícios, projetos, seminários, relatórios e estudos de caso, que poderão contribuir para a formação das notas ao longo do curso. A média final (MF) será calculada pela soma das notas obtidas (N), de acordo com a fórmula: (N1 + ... + Nn) / n.")
    ,@("History of microbiology, industrial microbiology, microbial phylogeny, characterization of microorganisms, nutrition and cultivation of microorganisms, viruses, filamentous fungi, yeasts, microalgae, bacteria.", "Provide students with fundamental knowledge of the main groups of microorganisms, covering their physiology, genetics, and roles in biotechnological processes, as well as developing practical skills in microbiological techniques for isolating, cultivating, and characterizing these organisms. Additionally, it encourages the critical interpretation of experimental results.")
    ,@("1. Histórico da microbiologia;`v2. As bases da microbiologia;`v3. Metabolismo microbiano;`v4. Crescimento microbiano;`v5. Controle do crescimento microbiano;`v6. Genética microbiana;`v7. Diversidade microbiana;`v8. Classificação dos microrganismos;`v9. Ecologia microbiana e microbiologia ambiental;`v10. Microbiologia industrial e aplicada.", "Média final (MF) ≥ 5,0 para aprovação. Prova de recuperação para alunos com MF<5,0")
    ,@("1. History of microbiology;`v2. Fundamentals of microbiology;`v3. Microbial metabolism;`v4. Microbial growth;`v5. Control of microbial growth;`v6. Microbial genetics;`v7 - Microbial diversity;`v8 - Classification of microorganisms;`v9 - Microbial ecology and environmental microbiology;`v10 - Industrial and applied microbiology.", "Theoretical Classes:`v1.History of microbiology, industrial importance of microorganisms, molecular taxonomy and phylogeny, physiology, and characterization of different groups: bacteria, fungi, viruses, and archaea;`v2.Microbial nutrition; culture media; environmental factors; microbial reproduction and growth; growth control measures; methods for microorganism quantification;`v3.Organization of the microbial genome. Genetic transfer: conjugation, transformation, and transduction. Gene expression regulation: Operon system;`v4.Ecological role of microorganisms in biocycles. Bioindicators; Symbiotic relationships;`v5.Viral genome and its diversity.`vPractical Classes:`v1.Optical microscopy and staining techniques;`v2.Preparation and sterilization of culture media - nutritional requirements;`v3.Microorganism inoculation techniques and pure culture isolation;`v4.Effects of physical and chemical agents on microbial growth;`v5.Microorganism quantification techniques: Colony-Forming Units (CFU), Turbidity, and Cell counting chamber under microscopy;`v6.Biochemical tests for microorganism identification and characterization;`v7.Bacterial analysis of water and milk samples.")
    ,@("A avaliação será feita por meio de provas escritas, trabalhos, seminários e participação.`v", "A média final após a recuperação será calculada como (MF + RC) / 2, sendo RC a nota da prova de recuperação. O aluno será aprovado se essa média final for ≥ 5,0.`v")
    ,@("A Nota final (NF) será calculada da seguinte maneira: NF = (P1 + P2)/2.`v", "1. PELCZAR Jr, M.J., CHAN, S.S., KRIEG, N.R. Microbiologia conceitos e aplicações, 2 ed. (Vol 1), São Paulo: Pearson Education do Brasil, 1997.`v2. MADIGAN, M.T., MARTINKO, J.M., PARKER, I. Microbiologia de Brock. São Paulo: Prentice Hall, 14a edição, 2016.`v3. TORTORA, G.J., FUNKE, B.R., CASE, C.L. Microbiologia, Artmed, Porto Alegre, RS, 12a edição, 2017.`v4. VERMELHO, A.B., PEREIRA, A., COELHO, R. e SOUTO-PADRON, T. Práticas de`vMicrobiologia, 2 edição. Rio de Janeiro: Guanabara-Koogan, 2019.`v")
    ,@("A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR).", "4873328 - Fernando Segato")
    ,@("1. PELCZAR Jr, M.J., CHAN, S.S., KRIEG, N.R. Microbiologia conceitos e aplicações, 2 ed. (Vol 1), São Paulo: Pearson Education do Brasil, 1997.`v2. MADIGAN, M.T., MARTINKO, J.M., PARKER, I. Microbiologia de Brock. São Paulo: Prentice Hall, 14a edição, 2016.`v3. TORTORA, G.J., FUNKE, B.R., CASE, C.L. Microbiologia, Artmed, Porto Alegre, RS, 12a edição, 2017.", "8853480 - Tatiane da Franca Silva")
)

foreach ($op in $ops) {
    $findText = $op[0]
    $replText = $op[1]
    $result = $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replText, 2)
    if (-not $result) {
        Write-Host "WARNING: replacement not found for: $findText"
    }
}

Write-Host "All replacements applied."
